$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1359.6066
$ws.Range("J17").Value = 1359.6066
$ws.Range("L17").Value = 4078.8198
$ws.Range("N17").Value = -4414.8198
$ws.Range("H19").Value = 1027.1818
$ws.Range("I19").Value = 950
$ws.Range("J19").Value = 1071.2858
$ws.Range("K19").Value = 950
$ws.Range("L19").Value = 1071.2858
$ws.Range("M19").Value = -775
$ws.Range("N19").Value = -1421.2858
$ws.Range("H40").Value = 16679571
$ws.Range("J40").Value = 33343162
$ws.Range("L40").Value = 33343162
$ws.Range("N40").Value = -33343512
$ws.Range("H43").Value = 3379522.5
$ws.Range("I43").Value = 5629204.5
$ws.Range("J43").Value = 4999.5
$ws.Range("K43").Value = 5629204.5
$ws.Range("L43").Value = 4999.5
$ws.Range("M43").Value = -5629135.5
$ws.Range("N43").Value = -5137.5
$ws.Range("H51").Value = 5531.0625
$ws.Range("I51").Value = 5135.7144
$ws.Range("K51").Value = 5135.7144
$ws.Range("M51").Value = -4651.7144
$ws.Range("H55").Value = 529.61536
$ws.Range("I55").Value = 710.4286
$ws.Range("J55").Value = 318.66666
$ws.Range("K55").Value = 710.4286
$ws.Range("L55").Value = 318.66666
$ws.Range("M55").Value = -496.4286
$ws.Range("N55").Value = -746.66666
$ws.Range("H62").Value = 3385.6428
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376
$ws.Range("H65").Value = 3385.6428
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880
$ws.Range("H86").Value = 2990569
$ws.Range("I86").Value = 6470412.5
$ws.Range("K86").Value = 6470412.5
$ws.Range("M86").Value = -6469289.5
$ws.Range("H89").Value = 2990569
$ws.Range("I89").Value = 6470412.5
$ws.Range("K89").Value = 32352062.5
$ws.Range("M89").Value = -32346446.5
$ws.Range("H98").Value = 2175.7058
$ws.Range("I98").Value = 1665.8
$ws.Range("K98").Value = 1665.8
$ws.Range("M98").Value = -167.8
$ws.Range("H106").Value = 41669910
$ws.Range("I106").Value = 47621332
$ws.Range("K106").Value = 47621332
$ws.Range("M106").Value = -47620701
$ws.Range("H111").Value = 4044.5625
$ws.Range("J111").Value = 4842.1816
$ws.Range("L111").Value = 14526.5448
$ws.Range("N111").Value = -20660.5448
$ws.Range("H122").Value = 2175.7058
$ws.Range("I122").Value = 1665.8
$ws.Range("K122").Value = 4997.4
$ws.Range("M122").Value = -2547.4
$ws.Range("H139").Value = 124195
$ws.Range("J139").Value = 124195
$ws.Range("L139").Value = 124195
$ws.Range("N139").Value = -134475
$ws.Range("H140").Value = 66827
$ws.Range("J140").Value = 66180
$ws.Range("L140").Value = 66180
$ws.Range("N140").Value = -76540

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13049.181
$ws.Range("I32").Value = 13247.407
$ws.Range("J32").Value = 12149.538
$ws.Range("K32").Value = 13247.407
$ws.Range("L32").Value = 12149.538
$ws.Range("M32").Value = -12960.407
$ws.Range("N32").Value = -12723.538
$ws.Range("H61").Value = 6904.683
$ws.Range("I61").Value = 7779
$ws.Range("K61").Value = 7779
$ws.Range("M61").Value = -7567
$ws.Range("H110").Value = 5545.5
$ws.Range("I110").Value = 3178.4546
$ws.Range("J110").Value = 10753
$ws.Range("K110").Value = 3178.4546
$ws.Range("L110").Value = 10753
$ws.Range("M110").Value = -1133.4546
$ws.Range("N110").Value = -14843
$ws.Range("H122").Value = 5983.1177
$ws.Range("I122").Value = 4892.091
$ws.Range("K122").Value = 14676.273
$ws.Range("M122").Value = -12226.273
$ws.Range("H132").Value = 15431.72
$ws.Range("I132").Value = 23326.54
$ws.Range("K132").Value = 69979.62
$ws.Range("M132").Value = -67449.62
$ws.Range("H136").Value = 6904.683
$ws.Range("I136").Value = 7779
$ws.Range("K136").Value = 23337
$ws.Range("M136").Value = -20787

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 25000
$ws.Range("J63").Value = 25000
$ws.Range("L63").Value = 25000
$ws.Range("N63").Value = -26372
$ws.Range("H66").Value = 25000
$ws.Range("J66").Value = 25000
$ws.Range("L66").Value = 75000
$ws.Range("N66").Value = -81864

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5765.2246
$ws.Range("I31").Value = 1517.4546
$ws.Range("J31").Value = 6994.8423
$ws.Range("K31").Value = 1517.4546
$ws.Range("L31").Value = 6994.8423
$ws.Range("M31").Value = -1222.4546
$ws.Range("N31").Value = -7584.8423
$ws.Range("H34").Value = 5765.2246
$ws.Range("I34").Value = 1517.4546
$ws.Range("J34").Value = 6994.8423
$ws.Range("K34").Value = 1517.4546
$ws.Range("L34").Value = 6994.8423
$ws.Range("M34").Value = -1315.4546
$ws.Range("N34").Value = -7398.8423
$ws.Range("H62").Value = 36937
$ws.Range("I62").Value = 3064.4
$ws.Range("K62").Value = 3064.4
$ws.Range("M62").Value = -2440.4
$ws.Range("H65").Value = 36937
$ws.Range("I65").Value = 3064.4
$ws.Range("K65").Value = 15322
$ws.Range("M65").Value = -12202
$ws.Range("H115").Value = 52999.25
$ws.Range("J115").Value = 52999.25
$ws.Range("L115").Value = 52999.25
$ws.Range("N115").Value = -55349.25
$ws.Range("H141").Value = 80866.60000000001
$ws.Range("J141").Value = 85165.766
$ws.Range("L141").Value = 85165.766
$ws.Range("N141").Value = -95525.766

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 315
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 315
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H139").Value = 2643.611
$ws.Range("I139").Value = 1780.6364
$ws.Range("K139").Value = 5341.9092
$ws.Range("M139").Value = -201.9092000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7559.375
$ws.Range("I126").Value = 5905.4546
$ws.Range("J126").Value = 11198
$ws.Range("K126").Value = 17716.3638
$ws.Range("L126").Value = 33594
$ws.Range("M126").Value = -15246.3638
$ws.Range("N126").Value = -38534

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7800.3
$ws.Range("I7").Value = 6000
$ws.Range("K7").Value = 6000
$ws.Range("M7").Value = -5888
$ws.Range("H93").Value = 995.2857
$ws.Range("I93").Value = 942.1579
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 942.1579
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = 305.8421
$ws.Range("N93").Value = -3996
$ws.Range("H122").Value = 41672468
$ws.Range("I122").Value = 71433170
$ws.Range("J122").Value = 7490
$ws.Range("K122").Value = 214299510
$ws.Range("L122").Value = 22470
$ws.Range("M122").Value = -214297060
$ws.Range("N122").Value = -27370
$ws.Range("H126").Value = 7800.3
$ws.Range("I126").Value = 6000
$ws.Range("K126").Value = 18000
$ws.Range("M126").Value = -15530
$ws.Range("H139").Value = 32825
$ws.Range("I139").Value = 32825
$ws.Range("K139").Value = 32825
$ws.Range("M139").Value = -27685

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I113").Value = 884.5
$ws.Range("J113").Value = 863.75
$ws.Range("K113").Value = 2653.5
$ws.Range("L113").Value = 2591.25
$ws.Range("M113").Value = -483.5
$ws.Range("N113").Value = -6931.25
$ws.Range("H132").Value = 7355365.5
$ws.Range("I132").Value = 794.56525
$ws.Range("K132").Value = 2383.69575
$ws.Range("M132").Value = 146.3042500000001
